$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose new value would otherwise be
# auto-converted to a number by Excel (e.g. "1.00", "0.0000237"), so the
# stored cell stays a literal text string identical to the source data.
$textCells = @(
    "D5", "D6", "D7", "D8", "D9", "D11", "D13", "D15", "D17", "D19",
    "D20", "D21", "D22", "D23", "D24", "D25", "D27", "D28", "D30", "D31",
    "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D41", "D43",
    "D44", "D45", "D46", "D47", "D48", "D49", "D50", "D51"
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "87.310.45"
$ws.Range("E2").Value = "  -1.38%  "

$ws.Range("D3").Value = "3.150.29"
$ws.Range("E3").Value = "  -7.27%  "

$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").Value = "205.85"
$ws.Range("E5").Value = "  -6.12%  "

$ws.Range("D6").Value = "605.93"
$ws.Range("E6").Value = "  -6.51%  "

$ws.Range("D7").Value = "0.371"
$ws.Range("E7").Value = "  -13.09%  "

$ws.Range("D8").Value = "0.657"
$ws.Range("E8").Value = "  -1.40%  "

$ws.Range("D9").Value = "1.00"
$ws.Range("E9").Value = "  -0.01%  "

$ws.Range("D10").Value = "3.146.27"
$ws.Range("E10").Value = "  -7.28%  "

$ws.Range("D11").Value = "0.529"
$ws.Range("E11").Value = "  -15.24%  "

$ws.Range("E12").Value = "  +4.24%  "

$ws.Range("D13").Value = "0.0000237"
$ws.Range("E13").Value = "  -17.92%  "

$ws.Range("D14").Value = "3.736.81"
$ws.Range("E14").Value = "  -7.29%  "

$ws.Range("D15").Value = "5.20"
$ws.Range("E15").Value = "  -7.02%  "

$ws.Range("D16").Value = "87.213.26"
$ws.Range("E16").Value = "  -1.38%  "

$ws.Range("D17").Value = "31.83"
$ws.Range("E17").Value = "  -13.42%  "

$ws.Range("D18").Value = "3.152.58"
$ws.Range("E18").Value = "  -7.33%  "

$ws.Range("D19").Value = "3.03"
$ws.Range("E19").Value = "  -0.90%  "

$ws.Range("D20").Value = "13.30"
$ws.Range("E20").Value = "  -11.41%  "

$ws.Range("D21").Value = "411.50"
$ws.Range("E21").Value = "  -10.83%  "

$ws.Range("D22").Value = "8.37"
$ws.Range("E22").Value = "  -14.05%  "

$ws.Range("D23").Value = "5.01"
$ws.Range("E23").Value = "  -11.22%  "

$ws.Range("D24").Value = "5.11"
$ws.Range("E24").Value = "  -7.95%  "

$ws.Range("D25").Value = "11.72"
$ws.Range("E25").Value = "  -7.49%  "

$ws.Range("D26").Value = "3.329.66"
$ws.Range("E26").Value = "  -7.08%  "

$ws.Range("D27").Value = "72.92"
$ws.Range("E27").Value = "  -10.37%  "

$ws.Range("D28").Value = "0.0000127"
$ws.Range("E28").Value = "  -11.82%  "

$ws.Range("E29").Value = "  +0.01%  "

$ws.Range("D30").Value = "1.00"
$ws.Range("E30").Value = "  -0.12%  "

$ws.Range("D31").Value = "0.157"
$ws.Range("E31").Value = "  -16.77%  "

$ws.Range("D32").Value = "538.82"
$ws.Range("E32").Value = "  -6.79%  "

$ws.Range("D33").Value = "8.11"
$ws.Range("E33").Value = "  -14.21%  "

$ws.Range("D34").Value = "1.28"
$ws.Range("E34").Value = "  -17.63%  "

$ws.Range("D35").Value = "1.83"
$ws.Range("E35").Value = "  -13.02%  "

$ws.Range("D36").Value = "6.62"
$ws.Range("E36").Value = "  -11.49%  "

$ws.Range("D37").Value = "0.131"
$ws.Range("E37").Value = "  -8.24%  "

$ws.Range("D38").Value = "21.59"
$ws.Range("E38").Value = "  -10.03%  "

$ws.Range("D39").Value = "21.83"
$ws.Range("E39").Value = "  -0.21%  "

$ws.Range("E40").Value = "  +0.04%  "

$ws.Range("D41").Value = "2.94"
$ws.Range("E41").Value = "  -6.22%  "

$ws.Range("D43").Value = "1.89"
$ws.Range("E43").Value = "  -11.04%  "

$ws.Range("D44").Value = "0.364"
$ws.Range("E44").Value = "  -16.94%  "

$ws.Range("D45").Value = "148.48"
$ws.Range("E45").Value = "  -5.83%  "

$ws.Range("D46").Value = "171.12"
$ws.Range("E46").Value = "  -9.71%  "

$ws.Range("D47").Value = "42.98"
$ws.Range("E47").Value = "  -8.10%  "

$ws.Range("D48").Value = "0.123"
$ws.Range("E48").Value = "  -1.59%  "

$ws.Range("D49").Value = "1.20"
$ws.Range("E49").Value = "  -15.92%  "

$ws.Range("D50").Value = "3.92"
$ws.Range("E50").Value = "  -13.49%  "

$ws.Range("B51").Value = "Mantle"
$ws.Range("C51").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D51").Value = "0.687"
$ws.Range("E51").Value = "  -13.00%  "

